# Season up to 1/17
# - "Games" sheet: the 1/16 game (vs SAC) has now been played; append it as
#   a new row (Game 40) at the bottom of the results table.
# - "Next" sheet: remove that now-played game from the top of the upcoming
#   schedule, shifting the remaining fixtures up by one row.

$wb = $excel.ActiveWorkbook

$games = $wb.Worksheets.Item("Games")
$next  = $wb.Worksheets.Item("Next")

# --- Games: append the result of the game that was played on 2024-01-16 ---
$newRow = $games.UsedRange.Rows.Count + 1

$games.Cells.Item($newRow, 1).Value  = 40        # Game
$games.Cells.Item($newRow, 2).Value  = 45307      # Date
$games.Cells.Item($newRow, 2).NumberFormat = $games.Cells.Item($newRow - 1, 2).NumberFormat
$games.Cells.Item($newRow, 3).Value  = 3          # Streak
$games.Cells.Item($newRow, 4).Value  = 119        # Pts
$games.Cells.Item($newRow, 5).Value  = 101.6      # Pace
$games.Cells.Item($newRow, 6).Value  = 0.594      # eFG
$games.Cells.Item($newRow, 7).Value  = 14.6       # TOV
$games.Cells.Item($newRow, 8).Value  = 23.5       # ORB
$games.Cells.Item($newRow, 9).Value  = 0.212      # FTR
$games.Cells.Item($newRow, 10).Value = 117.1      # ORT
$games.Cells.Item($newRow, 11).Value = "SAC"      # OppID
$games.Cells.Item($newRow, 12).Value = 117        # OppPts
$games.Cells.Item($newRow, 13).Value = 0.598      # OppeFG
$games.Cells.Item($newRow, 14).Value = 12.6       # OppTOV
$games.Cells.Item($newRow, 15).Value = 11.6       # OppORB
$games.Cells.Item($newRow, 16).Value = 0.076      # OppFTR
$games.Cells.Item($newRow, 17).Value = 115.1      # OppORT
$games.Cells.Item($newRow, 18).Value = 1          # Location
$games.Cells.Item($newRow, 19).Value = 1          # Target

# --- Next: drop the game that just got played (top row of the schedule) ---
$next.Rows.Item(2).Delete()
